# Regenerate save_data column G ("K") using K (strikeouts) instead of Strike# (pitch count),
# and write the recalculated values (s_vals) into the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G (K), rows 2-43, in order.
$newK = @(9,2,3,5,1,2,5,3,5,6,3,0,0,5,7,3,0,2,3,1,2,1,7,5,5,3,3,8,3,4,6,3,2,3,4,8,1,3,2,3,4,2)

$startRow = 2
for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
